$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format column A for the new rows as Text so the date-like strings
# ("2022/02/26" etc.) are stored as literal text, matching the existing rows,
# instead of being auto-converted into date serial numbers.
$ws.Range("A700:A709").NumberFormat = "@"

$ws.Cells.Item(700, 1).Value = "2022/02/26"
$ws.Cells.Item(700, 2).Value = 699
$ws.Cells.Item(700, 3).Value = 31201
$ws.Cells.Item(700, 4).Value = 642
$ws.Cells.Item(700, 5).Value = 73
$ws.Cells.Item(700, 6).Value = 0.02057626358
$ws.Cells.Item(700, 7).Value = 30418
$ws.Cells.Item(700, 8).Value = 27040
$ws.Cells.Item(700, 9).Value = 58241
$ws.Cells.Item(700, 10).Value = 0
$ws.Cells.Item(700, 11).Value = 0
$ws.Cells.Item(700, 12).Value = 0
$ws.Cells.Item(700, 13).Value = 571
$ws.Cells.Item(700, 14).Value = 0
$ws.Cells.Item(700, 15).Value = 571
$ws.Cells.Item(700, 16).Value = 8
$ws.Cells.Item(700, 17).Value = 6
$ws.Cells.Item(700, 18).Value = 11
$ws.Cells.Item(700, 19).Value = 56
$ws.Cells.Item(700, 20).Value = 100

$ws.Cells.Item(701, 1).Value = "2022/02/27"
$ws.Cells.Item(701, 2).Value = 700
$ws.Cells.Item(701, 3).Value = 31201
$ws.Cells.Item(701, 4).Value = 642
$ws.Cells.Item(701, 5).Value = 73
$ws.Cells.Item(701, 6).Value = 0.02057626358
$ws.Cells.Item(701, 7).Value = 30418
$ws.Cells.Item(701, 8).Value = 27040
$ws.Cells.Item(701, 9).Value = 58241
$ws.Cells.Item(701, 10).Value = 0
$ws.Cells.Item(701, 11).Value = 0
$ws.Cells.Item(701, 12).Value = 0
$ws.Cells.Item(701, 13).Value = 571
$ws.Cells.Item(701, 14).Value = 0
$ws.Cells.Item(701, 15).Value = 571
$ws.Cells.Item(701, 16).Value = 7
$ws.Cells.Item(701, 17).Value = 6
$ws.Cells.Item(701, 18).Value = 11
$ws.Cells.Item(701, 19).Value = 56
$ws.Cells.Item(701, 20).Value = 101

$ws.Cells.Item(702, 1).Value = "2022/02/28"
$ws.Cells.Item(702, 2).Value = 701
$ws.Cells.Item(702, 3).Value = 31201
$ws.Cells.Item(702, 4).Value = 642
$ws.Cells.Item(702, 5).Value = 73
$ws.Cells.Item(702, 6).Value = 0.02057626358
$ws.Cells.Item(702, 7).Value = 30418
$ws.Cells.Item(702, 8).Value = 27040
$ws.Cells.Item(702, 9).Value = 58241
$ws.Cells.Item(702, 10).Value = 0
$ws.Cells.Item(702, 11).Value = 0
$ws.Cells.Item(702, 12).Value = 0
$ws.Cells.Item(702, 13).Value = 571
$ws.Cells.Item(702, 14).Value = 0
$ws.Cells.Item(702, 15).Value = 571
$ws.Cells.Item(702, 16).Value = 8
$ws.Cells.Item(702, 17).Value = 6
$ws.Cells.Item(702, 18).Value = 11
$ws.Cells.Item(702, 19).Value = 56
$ws.Cells.Item(702, 20).Value = 101

$ws.Cells.Item(703, 1).Value = "2022/03/01"
$ws.Cells.Item(703, 2).Value = 702
$ws.Cells.Item(703, 3).Value = 31201
$ws.Cells.Item(703, 4).Value = 68
$ws.Cells.Item(703, 5).Value = 73
$ws.Cells.Item(703, 6).Value = 0.002179417326
$ws.Cells.Item(703, 7).Value = 30418
$ws.Cells.Item(703, 8).Value = 27040
$ws.Cells.Item(703, 9).Value = 58241
$ws.Cells.Item(703, 10).Value = 0
$ws.Cells.Item(703, 11).Value = -574
$ws.Cells.Item(703, 12).Value = 0
$ws.Cells.Item(703, 13).Value = 571
$ws.Cells.Item(703, 14).Value = 0
$ws.Cells.Item(703, 15).Value = 571
$ws.Cells.Item(703, 16).Value = 8
$ws.Cells.Item(703, 17).Value = 6
$ws.Cells.Item(703, 18).Value = 11
$ws.Cells.Item(703, 19).Value = 56
$ws.Cells.Item(703, 20).Value = 101

$ws.Cells.Item(704, 1).Value = "2022/03/02"
$ws.Cells.Item(704, 2).Value = 703
$ws.Cells.Item(704, 3).Value = 31201
$ws.Cells.Item(704, 4).Value = 642
$ws.Cells.Item(704, 5).Value = 73
$ws.Cells.Item(704, 6).Value = 0.02057626358
$ws.Cells.Item(704, 7).Value = 30418
$ws.Cells.Item(704, 8).Value = 27040
$ws.Cells.Item(704, 9).Value = 58241
$ws.Cells.Item(704, 10).Value = 0
$ws.Cells.Item(704, 11).Value = 574
$ws.Cells.Item(704, 12).Value = 0
$ws.Cells.Item(704, 13).Value = 571
$ws.Cells.Item(704, 14).Value = 0
$ws.Cells.Item(704, 15).Value = 571
$ws.Cells.Item(704, 16).Value = 8
$ws.Cells.Item(704, 17).Value = 6
$ws.Cells.Item(704, 18).Value = 11
$ws.Cells.Item(704, 19).Value = 56
$ws.Cells.Item(704, 20).Value = 101

$ws.Cells.Item(705, 1).Value = "2022/03/03"
$ws.Cells.Item(705, 2).Value = 704
$ws.Cells.Item(705, 3).Value = 31613
$ws.Cells.Item(705, 4).Value = 644
$ws.Cells.Item(705, 5).Value = 8
$ws.Cells.Item(705, 6).Value = 0.02037136621
$ws.Cells.Item(705, 7).Value = 30893
$ws.Cells.Item(705, 8).Value = 27040
$ws.Cells.Item(705, 9).Value = 58653
$ws.Cells.Item(705, 10).Value = 412
$ws.Cells.Item(705, 11).Value = 2
$ws.Cells.Item(705, 12).Value = 412
$ws.Cells.Item(705, 13).Value = 373
$ws.Cells.Item(705, 14).Value = 0
$ws.Cells.Item(705, 15).Value = 373
$ws.Cells.Item(705, 16).Value = 4
$ws.Cells.Item(705, 17).Value = 5
$ws.Cells.Item(705, 18).Value = 11
$ws.Cells.Item(705, 19).Value = 56
$ws.Cells.Item(705, 20).Value = 101

$ws.Cells.Item(706, 1).Value = "2022/03/04"
$ws.Cells.Item(706, 2).Value = 705
$ws.Cells.Item(706, 3).Value = 31913
$ws.Cells.Item(706, 4).Value = 646
$ws.Cells.Item(706, 5).Value = 11
$ws.Cells.Item(706, 6).Value = 0.02024253439
$ws.Cells.Item(706, 7).Value = 31188
$ws.Cells.Item(706, 8).Value = 27040
$ws.Cells.Item(706, 9).Value = 58953
$ws.Cells.Item(706, 10).Value = 300
$ws.Cells.Item(706, 11).Value = 2
$ws.Cells.Item(706, 12).Value = 300
$ws.Cells.Item(706, 13).Value = 288
$ws.Cells.Item(706, 14).Value = 0
$ws.Cells.Item(706, 15).Value = 288
$ws.Cells.Item(706, 16).Value = 4
$ws.Cells.Item(706, 17).Value = 3
$ws.Cells.Item(706, 18).Value = 11
$ws.Cells.Item(706, 19).Value = 56
$ws.Cells.Item(706, 20).Value = 101

$ws.Cells.Item(707, 1).Value = "2022/03/05"
$ws.Cells.Item(707, 2).Value = 706
$ws.Cells.Item(707, 3).Value = 31913
$ws.Cells.Item(707, 4).Value = 646
$ws.Cells.Item(707, 5).Value = 11
$ws.Cells.Item(707, 6).Value = 0.02024253439
$ws.Cells.Item(707, 7).Value = 31188
$ws.Cells.Item(707, 8).Value = 27040
$ws.Cells.Item(707, 9).Value = 58953
$ws.Cells.Item(707, 10).Value = 0
$ws.Cells.Item(707, 11).Value = 0
$ws.Cells.Item(707, 12).Value = 0
$ws.Cells.Item(707, 13).Value = 288
$ws.Cells.Item(707, 14).Value = 0
$ws.Cells.Item(707, 15).Value = 288
$ws.Cells.Item(707, 16).Value = 4
$ws.Cells.Item(707, 17).Value = 3
$ws.Cells.Item(707, 18).Value = 11
$ws.Cells.Item(707, 19).Value = 56
$ws.Cells.Item(707, 20).Value = 101

$ws.Cells.Item(708, 1).Value = "2022/03/06"
$ws.Cells.Item(708, 2).Value = 707
$ws.Cells.Item(708, 3).Value = 31913
$ws.Cells.Item(708, 4).Value = 646
$ws.Cells.Item(708, 5).Value = 11
$ws.Cells.Item(708, 6).Value = 0.02024253439
$ws.Cells.Item(708, 7).Value = 31188
$ws.Cells.Item(708, 8).Value = 27040
$ws.Cells.Item(708, 9).Value = 58953
$ws.Cells.Item(708, 10).Value = 0
$ws.Cells.Item(708, 11).Value = 0
$ws.Cells.Item(708, 12).Value = 0
$ws.Cells.Item(708, 13).Value = 288
$ws.Cells.Item(708, 14).Value = 0
$ws.Cells.Item(708, 15).Value = 288
$ws.Cells.Item(708, 16).Value = 4
$ws.Cells.Item(708, 17).Value = 3
$ws.Cells.Item(708, 18).Value = 11
$ws.Cells.Item(708, 19).Value = 56
$ws.Cells.Item(708, 20).Value = 102

$ws.Cells.Item(709, 1).Value = "2022/03/07"
$ws.Cells.Item(709, 2).Value = 708
$ws.Cells.Item(709, 3).Value = 32278
$ws.Cells.Item(709, 4).Value = 647
$ws.Cells.Item(709, 5).Value = 6
$ws.Cells.Item(709, 6).Value = 0.02004461243
$ws.Cells.Item(709, 7).Value = 31557
$ws.Cells.Item(709, 8).Value = 27040
$ws.Cells.Item(709, 9).Value = 59318
$ws.Cells.Item(709, 10).Value = 365
$ws.Cells.Item(709, 11).Value = 1
$ws.Cells.Item(709, 12).Value = 365
$ws.Cells.Item(709, 13).Value = 275
$ws.Cells.Item(709, 14).Value = 0
$ws.Cells.Item(709, 15).Value = 275
$ws.Cells.Item(709, 16).Value = 10
$ws.Cells.Item(709, 17).Value = 1
$ws.Cells.Item(709, 18).Value = 11
$ws.Cells.Item(709, 19).Value = 56
$ws.Cells.Item(709, 20).Value = 102

# Reset the number format back to General/Normal now that the text values
# are committed, so the cells carry no explicit style (matching the target).
$ws.Range("A700:A709").Style = "Normal"
